# Regenerate the "K" column (column G) values for rows 2-57 on the active
# sheet. This corresponds to "regen save_data to use K instead of Strike#,
# regen std/mean, calc and write s_vals" - the K column values were
# recomputed and rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G (header "K"), for rows 2 through 57, in order.
$kValues = @(1, 1, 2, 1, 1, 1, 0, 1, 1, 0, 0, 1, 1, 1, 1, 0, 2, 1, 2, 0, 2, 5, 2, 1, 2, 2, 0, 0, 3, 1, 1, 0, 1, 3, 0, 0, 1, 2, 0, 4, 3, 0, 0, 2, 0, 2, 1, 3, 0, 3, 1, 2, 1, 1, 0, 2)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
